$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$newValues = @(
    "372cb3aa-d3cc-455b-a905-55d42e516ed1",
    "0ed122a5-5c70-4cad-9cec-53362f0d714b",
    "0637a275-a245-4bd1-904d-38462b3f8c12",
    "c812b8c4-786a-4566-84b2-207142fe1175",
    "85c0de35-8816-4631-abfb-2d3a290d3d59",
    "a2ae6296-de2b-45e3-862a-a1fa81caed99",
    "83c2aa69-9cd6-4379-bf1c-58680e3a79e1",
    "0041ec7b-0505-4fd9-9cf1-98bab786f898",
    "4372dc0b-137f-425e-862b-8779645a1edd",
    "2daf847a-7c6c-4430-a71d-ec28c5fea551",
    "9bdc4625-b04c-4ff8-b01f-cf79aa0adf7d",
    "701dc41b-d2c1-4d97-bf1c-d348b334904b",
    "2159dc03-98c8-4bd5-ae57-525d75881eba",
    "ae3bbbe0-944c-44e7-926a-50915cfc6e1f",
    "3a591750-d1ce-4d62-9669-efbf6ee90631",
    "db8aabcb-8085-4dbd-8e49-ebd9364c181f",
    "ed0a18d5-11cd-4fc4-ba5d-adf873b7595c",
    "b96fd29f-efae-45c4-ac8a-80420a6c5936",
    "cd278cf5-384d-4729-996c-0420f5d9d628",
    "9a3a41b1-2a42-452a-8c27-b29fe55214f4",
    "d70754be-cc55-41ec-b06d-96243a529b1e",
    "24a5bcc3-ddb9-4a9e-b730-0bebd0b42617",
    "73706993-b444-4171-bca6-566335e3fffc",
    "0046b82e-9ab4-43c5-bf80-ecd54c160e30",
    "bab249b5-51eb-444b-818b-de6dc059ff41",
    "d92f112a-93d1-4bdd-bc82-ac2f0d8d5d2e",
    "6a4ec714-33e6-44e3-a238-0e60838c8a8c",
    "59aa3cc6-efe0-40b3-b430-92f4278e0a3a",
    "d0197805-5078-40ff-82da-52ab10602dbf",
    "975f06e4-d477-4015-9a9d-9ddf165ed0c9",
    "b1fcdaf0-bcb4-464b-8e6d-b8e2abb5ccc8",
    "0e5864ae-9696-49dd-866b-1224bb8ca615",
    "2907bfbd-9c8b-45d4-bc74-12fb54d3fd24",
    "b4a27f1e-19fe-4c6b-8ce8-3e60f9d85655",
    "9f36a0c8-3e4e-4ee5-9a83-a905decfa0d8",
    "82c3873f-b121-42e2-9bf0-9806370453b3",
    "72929b93-50b1-4489-8a67-73bb09af0ff2",
    "c72075d5-a03f-4a50-9f7c-3674efc50f8a",
    "b084fa4d-c34b-4e49-bef4-6088df401a96",
    "763b57a5-c33c-442b-88d3-f2a9b1a99213",
    "ac910329-0262-4d62-84a8-d2c6bd40638e",
    "91b551a8-6f96-433e-bd34-4a8744009d83",
    "dc749e63-1e6d-44b8-9c79-361abed985a7",
    "ad3004cb-1002-40b9-baca-169df14c1700",
    "c94524a8-53a1-4f75-8aba-a3e137cef2e6",
    "d746cfb4-4bfc-4dfd-9d8d-eaa011bec22b",
    "49899fc0-ad1a-4a1e-b72d-0b9e1888b952",
    "378a9579-bd6a-44d0-b9cb-bdcd2ff88270",
    "ecdd500f-7b06-40da-b868-c3bb5ca413d4",
    "0f9eaa0e-6b7e-4b54-a464-27abcd79d613",
    "c9e048e4-389e-4274-b9c4-11e86d7337ee",
    "2ee0b975-83d1-4fe6-9a2a-7cabbbd15cb9",
    "8f62aacd-7797-482e-bf40-adebd0532de2",
    "c4e64aa2-49a1-4f3d-9074-90f7f756e3fd",
    "0b981b84-11ed-4c66-b8c5-75626c99d5e7",
    "d735378e-eb64-42a4-80f4-fc4e175fa6aa",
    "4a31bde2-3396-4e79-92dd-6cd9f3a9f21a",
    "025cabee-24af-4ff9-99a8-31c9c1adde43",
    "554e4017-bd12-4917-847c-b755f3b47824",
    "f90c2bef-ee88-44e6-b475-f9cff7f97a65",
    "d4ce1835-46f4-49c7-a8a9-7447a67b6b16",
    "693431ed-8f29-4f7f-8906-48b117729456",
    "9f8d9646-eaf1-4ec4-8434-4a2c1a70ec40",
    "e26a68bc-aa42-4fc5-9bc8-d669f2d8366a",
    "a62f8413-a713-4c95-9a7e-2093a36e58ad",
    "4da7cbcb-8083-483a-be7c-67c7a6baac6a",
    "70c885bf-80a6-42d8-b0d4-c8897471757e",
    "4929f468-ab9f-46db-81f7-c5c368329ce0",
    "2fcf07f2-20fd-47c8-8218-6a92d1109299",
    "37af9c50-b101-4988-92bc-7c32844139d1",
    "1f4c623b-b4cd-439a-ab5c-43b4cefae2d3",
    "09da1850-3489-471e-aa5a-b062ba695ff3",
    "56cf1db1-c906-4a84-a31f-d3cd405132bc",
    "f4c55d4e-ddcc-4948-9bdb-772896b2630f",
    "71f17a13-d9a6-46fe-932f-0d0cfbf4cd66",
    "97b346ef-5a26-42a9-be25-4820c856a724",
    "2076b7b3-18b8-401a-a8ce-df6031bec2d8",
    "bbc9b991-5063-4394-bd46-c4cce43d9bcb",
    "54e0d1e2-b526-4ff8-b88b-1bb3eb75947b",
    "efaaa067-0e50-47e6-acd9-0ab3a5ce05fc",
    "3fdab6c9-39a3-4905-9019-33f3dd3711db",
    "c2b2ec59-31ba-4ab4-a006-2c95b3cb0dd8",
    "0593ced6-b74c-4d75-a958-59f034317205",
    "419a5c1d-c836-4e6e-9870-97c855ea29fa",
    "3f7f17e7-56e7-4b1d-be14-bb1e9239d956",
    "3f89b6b0-4b5d-49f7-9b11-81f91f9511ca",
    "b47ecb57-539e-4fcb-8fbc-a1a1fe0c6c67",
    "8d5a54ef-5289-4c7f-80a0-48359746709f",
    "3ab0389a-057e-4eef-9135-78fa900f5810",
    "83257be3-215b-4065-a54a-0ad6682cf23a",
    "32a94177-aa6e-4e52-bb40-4da5966c144c",
    "6568617e-58a7-405d-b583-ed4e86551044",
    "1ac97b52-297e-4ecb-bdac-c34828e7d60c",
    "7acdac34-7ae3-4e42-95b9-dc3ec74e5f1a",
    "b2c22bec-fa9a-4ed8-a546-de2835af7776",
    "67066532-86a1-4645-87d5-0ce62b34d801",
    "2a75fb98-d523-47ba-9141-598845efde37",
    "79575c63-9e16-41db-8e95-febc68cdb703",
    "f6ff5e47-86a5-4afa-a442-b19a87d2f3b9",
    "e5b5e157-8fe5-43c0-9649-6ccd9af6aa44",
    "3a9eb28a-f9fc-4b1a-a3c5-7eacfe0f5188",
    "9b6d9dcd-e9ad-45cb-90c5-88c941a832c7",
    "770bccb5-9e72-4e5f-b2ed-7016cd145c75",
    "267f305b-6cd0-40af-bb69-76c45cd358e2",
    "73af36cd-a77a-4e3d-8840-503958201bc0",
    "619d6334-c993-408c-a7a6-1cf7ed1db415",
    "5c3cad25-d713-4d4f-9ca3-fc8275069a71",
    "6c8ff5fa-9e7d-4de7-b381-1b66b81672bb",
    "f3cf86ad-0b3a-43a9-8733-d548e2094dad",
    "ad9bfcd1-e193-496e-8756-27c8d38ea84a",
    "b110e34b-b1a5-46ee-9971-7d4535561f28",
    "9df52aea-f06e-403a-a82e-9010d3adc9d2",
    "2597f3e4-7cd9-44ff-8159-2da33dfbfcc9",
    "c97455cd-54b3-495f-b762-10a1a960415f",
    "1885623d-c36a-41f7-8fed-9f303a0897b1",
    "c0ea0aa3-962f-4869-aa6b-fc368a8cc634",
    "34a0f612-b3f8-4888-86ce-28ca5ba496fd"
)

$startRow = 21
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 11).Value = $newValues[$i]
}
